# Generate Report for Handoff
# Refresh the "Latest Handoff Datetime" (column D, row 5 - the
# 3521f6a1-... file) on both the zh-cn and de-de status sheets with the
# timestamps of the new handoff that was just generated.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-04 07:39:53"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-04 07:40:09"
